# Update the date line, then update the 25 division-problem cells in the
# single 20-row x 5-col table. Only the first five problem-rows (1, 5, 9,
# 13, 17) carry text; the rows in between are blank answer rows.
#
# Cell.Range.Text is used (instead of a blanket Find/Replace) because the
# old/new problem strings are not unique across the table (e.g. "48÷5=9, 3"
# and "56÷9=6, 2" each show up as both a source value in one cell and a
# target value in another), so addressing must be done by cell position.

$d = $word.ActiveDocument

# 1) Header date.
$d.Content.Find.Execute("2025-04-05 Saturday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2025-04-06 Sunday", 2)

# 2) Division problems, in row-major order matching the diff.
$t = $d.Tables.Item(1)

$values = @(
    "14÷5=2, 4",
    "48÷5=9, 3",
    "59÷2=29, 1",
    "65÷6=10, 5",
    "79÷5=15, 4",
    "47÷7=6, 5",
    "36÷9=4, 0",
    "16÷2=8, 0",
    "37÷7=5, 2",
    "74÷2=37, 0",
    "40÷4=10, 0",
    "31÷6=5, 1",
    "41÷5=8, 1",
    "23÷4=5, 3",
    "93÷7=13, 2",
    "96÷9=10, 6",
    "51÷3=17, 0",
    "46÷4=11, 2",
    "22÷4=5, 2",
    "23÷4=5, 3",
    "56÷9=6, 2",
    "17÷8=2, 1",
    "16÷4=4, 0",
    "39÷5=7, 4",
    "90÷8=11, 2"
)

$rows = @(1, 5, 9, 13, 17)
$i = 0
foreach ($r in $rows) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$i]
        $i++
    }
}

Write-Output "updated date + $i problem cells"
